# Formed the consolidated report
# Set the "Absent" (column H) values for the rows that were still using
# placeholder/blank values, to complete the consolidated attendance report.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H3").Value  = 1
$ws.Range("H6").Value  = 0
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("H16").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 0
